$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1 = Wins, AE1 = Losses, AF1 = Ties, styled like the other
# header cells (bold font, border, centered) by copying the existing
# header cell's format instead of inventing a new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-54: record columns (Wins, Losses, Ties) as numbers.
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 72   # AD
    $ws.Cells.Item($r, 31).Value = 90   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
